# Fill in the confusion-matrix proportions + per-class/per-fold summary values
# that were still blank placeholders in the t-SNE results workbook.
$wb = $excel.ActiveWorkbook

# Keep calculation manual while we write data so the pre-existing stale
# #DIV/0! caches on the F/M ratio formulas are left untouched, matching
# the authored diff (those formula cells are not part of this edit).
$excel.Calculation = -4135  # xlCalculationManual

# --- Sheet: FS ---
$ws = $wb.Worksheets.Item("FS")
$ws.Range("C24").Value = 0.6
$ws.Range("D24").Value = 0.4
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 26.31578947368421
$ws.Range("J24").Value = 0.9777777777777777
$ws.Range("K24").Value = 0.022222222222222223
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 4.093567251461988
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.875
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.027777777777777776
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.05555555555555555
$ws.Range("L26").Value = 0.9259259259259259

# --- Sheet: IF ---
$ws = $wb.Worksheets.Item("IF")
$ws.Range("C18").Value = 1.0
$ws.Range("D18").Value = 0.0
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 36.84210526315789
$ws.Range("J18").Value = 0.8888888888888888
$ws.Range("K18").Value = 0.06666666666666667
$ws.Range("L18").Value = 0.044444444444444446
$ws.Range("N18").Value = 14.619883040935672
$ws.Range("C19").Value = 0.25
$ws.Range("D19").Value = 0.625
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.1111111111111111
$ws.Range("K19").Value = 0.8333333333333334
$ws.Range("L19").Value = 0.05555555555555555
$ws.Range("C20").Value = 0.3333333333333333
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.3333333333333333
$ws.Range("J20").Value = 0.037037037037037035
$ws.Range("K20").Value = 0.1111111111111111
$ws.Range("L20").Value = 0.8518518518518519
$ws.Range("C24").Value = 0.4
$ws.Range("D24").Value = 0.4
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 31.57894736842105
$ws.Range("J24").Value = 0.8666666666666667
$ws.Range("K24").Value = 0.13333333333333333
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 10.526315789473683
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.75
$ws.Range("E25").Value = 0.25
$ws.Range("J25").Value = 0.041666666666666664
$ws.Range("K25").Value = 0.9166666666666666
$ws.Range("L25").Value = 0.041666666666666664
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.8333333333333334
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.07407407407407407
$ws.Range("L26").Value = 0.8888888888888888

# --- Sheet: IA ---
$ws = $wb.Worksheets.Item("IA")
$ws.Range("C18").Value = 0.6
$ws.Range("D18").Value = 0.4
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 15.789473684210526
$ws.Range("J18").Value = 0.9777777777777777
$ws.Range("K18").Value = 0.022222222222222223
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 2.923976608187134
$ws.Range("C19").Value = 0.0
$ws.Range("D19").Value = 1.0
$ws.Range("E19").Value = 0.0
$ws.Range("J19").Value = 0.027777777777777776
$ws.Range("K19").Value = 0.9583333333333334
$ws.Range("L19").Value = 0.013888888888888888
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("E20").Value = 0.8333333333333334
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.018518518518518517
$ws.Range("L20").Value = 0.9814814814814815
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.9555555555555556
$ws.Range("K24").Value = 0.022222222222222223
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 4.093567251461988
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.027777777777777776
$ws.Range("C26").Value = 0.3333333333333333
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.3333333333333333
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 0.05555555555555555
$ws.Range("L26").Value = 0.9444444444444444

# --- Sheet: FS-IF ---
$ws = $wb.Worksheets.Item("FS-IF")
$ws.Range("C18").Value = 0.8
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 26.31578947368421
$ws.Range("J18").Value = 0.9777777777777777
$ws.Range("K18").Value = 0.0
$ws.Range("L18").Value = 0.022222222222222223
$ws.Range("N18").Value = 2.3391812865497075
$ws.Range("C19").Value = 0.125
$ws.Range("D19").Value = 0.625
$ws.Range("E19").Value = 0.25
$ws.Range("J19").Value = 0.027777777777777776
$ws.Range("K19").Value = 0.9722222222222222
$ws.Range("L19").Value = 0.0
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("E20").Value = 0.8333333333333334
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.018518518518518517
$ws.Range("L20").Value = 0.9814814814814815
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 36.84210526315789
$ws.Range("J24").Value = 0.9555555555555556
$ws.Range("K24").Value = 0.022222222222222223
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 5.847953216374268
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.75
$ws.Range("E25").Value = 0.25
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.027777777777777776
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.5
$ws.Range("E26").Value = 0.3333333333333333
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.07407407407407407
$ws.Range("L26").Value = 0.8888888888888888

# --- Sheet: FS-IA ---
$ws = $wb.Worksheets.Item("FS-IA")
$ws.Range("C18").Value = 0.8
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 21.052631578947366
$ws.Range("J18").Value = 0.9555555555555556
$ws.Range("K18").Value = 0.044444444444444446
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 1.1695906432748537
$ws.Range("C19").Value = 0.125
$ws.Range("D19").Value = 0.75
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.0
$ws.Range("K19").Value = 1.0
$ws.Range("L19").Value = 0.0
$ws.Range("C20").Value = 0.16666666666666666
$ws.Range("D20").Value = 0.0
$ws.Range("E20").Value = 0.8333333333333334
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.0
$ws.Range("L20").Value = 1.0
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.9777777777777777
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 5.263157894736842
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9583333333333334
$ws.Range("L25").Value = 0.041666666666666664
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.5
$ws.Range("E26").Value = 0.3333333333333333
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.07407407407407407
$ws.Range("L26").Value = 0.9074074074074074

# --- Sheet: IF-IA ---
$ws = $wb.Worksheets.Item("IF-IA")
$ws.Range("C18").Value = 1.0
$ws.Range("D18").Value = 0.0
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 21.052631578947366
$ws.Range("J18").Value = 1.0
$ws.Range("K18").Value = 0.0
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 0.5847953216374269
$ws.Range("C19").Value = 0.25
$ws.Range("D19").Value = 0.625
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.0
$ws.Range("K19").Value = 1.0
$ws.Range("L19").Value = 0.0
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("E20").Value = 0.8333333333333334
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.018518518518518517
$ws.Range("L20").Value = 0.9814814814814815
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 36.84210526315789
$ws.Range("J24").Value = 0.9777777777777777
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 3.508771929824561
$ws.Range("C25").Value = 0.25
$ws.Range("D25").Value = 0.375
$ws.Range("E25").Value = 0.375
$ws.Range("J25").Value = 0.013888888888888888
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.013888888888888888
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.8333333333333334
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.018518518518518517
$ws.Range("L26").Value = 0.9444444444444444

# --- Sheet: FS-IF-IA ---
$ws = $wb.Worksheets.Item("FS-IF-IA")
$ws.Range("C18").Value = 1.0
$ws.Range("D18").Value = 0.0
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 21.052631578947366
$ws.Range("J18").Value = 0.9777777777777777
$ws.Range("K18").Value = 0.022222222222222223
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 0.5847953216374269
$ws.Range("C19").Value = 0.125
$ws.Range("D19").Value = 0.75
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.0
$ws.Range("K19").Value = 1.0
$ws.Range("L19").Value = 0.0
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.6666666666666666
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.0
$ws.Range("L20").Value = 1.0
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 1.0
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 3.508771929824561
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.027777777777777776
$ws.Range("C26").Value = 0.3333333333333333
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.5
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.05555555555555555
$ws.Range("L26").Value = 0.9259259259259259

